$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.403.80"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.498.66"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.14%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.212"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.650"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000303"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.58"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.069.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "610.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.520.44"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.66"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.504.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.87%  "

$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.985"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.63"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.69"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.115"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.42"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +12.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.10"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.22%  "

$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "525.40"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.393"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.57"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.537.33"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.141"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.40"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0764"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0456"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.91"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.142"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.85"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.22%  "

$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.99"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.34"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.57%  "
